$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 9).Value = 0.9371385197543067
$ws.Cells.Item(2, 10).Value = 0.9371385197543067
$ws.Cells.Item(2, 13).Value = 1.845768666666667
$ws.Cells.Item(2, 14).Value = 5.537306
$ws.Cells.Item(2, 15).Value = 0.01459089321241885
$ws.Cells.Item(2, 16).Value = 0.01459089321241885
$ws.Cells.Item(2, 17).Value = 25.46532706448355
$ws.Cells.Item(2, 18).Value = 229.187943580352
$ws.Cells.Item(2, 19).Value = 0.01367368806697936
$ws.Cells.Item(2, 20).Value = 0.01367368806697936
$ws.Cells.Item(3, 9).Value = 0.9371385197543067
$ws.Cells.Item(3, 10).Value = 0.9371385197543067
$ws.Cells.Item(3, 15).Value = 0.6557810310272387
$ws.Cells.Item(3, 16).Value = 0.6557810310272387
$ws.Cells.Item(3, 19).Value = 0.6145576646998195
$ws.Cells.Item(3, 20).Value = 0.6145576646998195
$ws.Cells.Item(4, 9).Value = 0.9371385197543067
$ws.Cells.Item(4, 10).Value = 0.9371385197543067
$ws.Cells.Item(4, 13).Value = 41.69841866666667
$ws.Cells.Item(4, 14).Value = 125.095256
$ws.Cells.Item(4, 15).Value = 0.3296280757603424
$ws.Cells.Item(4, 16).Value = 0.3296280757603424
$ws.Cells.Item(4, 17).Value = 575.2962917807503
$ws.Cells.Item(4, 18).Value = 5177.666626026752
$ws.Cells.Item(4, 19).Value = 0.3089071669875078
$ws.Cells.Item(4, 20).Value = 0.3089071669875078
$ws.Cells.Item(5, 7).Value = 0.9254496666666667
$ws.Cells.Item(5, 8).Value = 2.776349
$ws.Cells.Item(5, 9).Value = 0.0628614802456932
$ws.Cells.Item(5, 10).Value = 0.06286148024569319
$ws.Cells.Item(5, 13).Value = 1.845768666666667
$ws.Cells.Item(5, 14).Value = 5.537306
$ws.Cells.Item(5, 15).Value = 0.01459089321241885
$ws.Cells.Item(5, 16).Value = 0.01459089321241885
$ws.Cells.Item(5, 17).Value = 1.708165997310445
$ws.Cells.Item(5, 18).Value = 15.373493975794
$ws.Cells.Item(5, 19).Value = 0.0009172051454394866
$ws.Cells.Item(5, 20).Value = 0.0009172051454394863
$ws.Cells.Item(6, 7).Value = 0.9254496666666667
$ws.Cells.Item(6, 8).Value = 2.776349
$ws.Cells.Item(6, 9).Value = 0.0628614802456932
$ws.Cells.Item(6, 10).Value = 0.06286148024569319
$ws.Cells.Item(6, 15).Value = 0.6557810310272387
$ws.Cells.Item(6, 16).Value = 0.6557810310272387
$ws.Cells.Item(6, 17).Value = 76.77274054260678
$ws.Cells.Item(6, 18).Value = 690.9546648834611
$ws.Cells.Item(6, 19).Value = 0.04122336632741909
$ws.Cells.Item(6, 20).Value = 0.04122336632741908
$ws.Cells.Item(7, 7).Value = 0.9254496666666667
$ws.Cells.Item(7, 8).Value = 2.776349
$ws.Cells.Item(7, 9).Value = 0.0628614802456932
$ws.Cells.Item(7, 10).Value = 0.06286148024569319
$ws.Cells.Item(7, 13).Value = 41.69841866666667
$ws.Cells.Item(7, 14).Value = 125.095256
$ws.Cells.Item(7, 15).Value = 0.3296280757603424
$ws.Cells.Item(7, 16).Value = 0.3296280757603424
$ws.Cells.Item(7, 17).Value = 38.58978765559378
$ws.Cells.Item(7, 18).Value = 347.308088900344
$ws.Cells.Item(7, 19).Value = 0.02072090877283463
$ws.Cells.Item(7, 20).Value = 0.02072090877283462
